$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.208955223880597
$ws.Range("C2").Value = 0.5335820895522388
$ws.Range("J2").Value = 0.01865671641791045
$ws.Range("P2").Value = 0.1492537313432836
$ws.Range("S2").Value = 0.08955223880597014
$ws.Range("B3").Value = 0.006756756756756757
$ws.Range("C3").Value = 0.02027027027027027
$ws.Range("J3").Value = 0.0472972972972973
$ws.Range("P3").Value = 0.7432432432432432
$ws.Range("S3").Value = 0.1824324324324324
$ws.Range("B6").Value = 0.07692307692307693
$ws.Range("D6").Value = 0.01282051282051282
$ws.Range("F6").Value = 0.1111111111111111
$ws.Range("J6").Value = 0.217948717948718
$ws.Range("O6").Value = 0.01282051282051282
$ws.Range("Q6").Value = 0.09829059829059829
$ws.Range("R6").Value = 0.06837606837606838
$ws.Range("S6").Value = 0.4017094017094017
$ws.Range("B7").Value = 0.1212121212121212
$ws.Range("D7").Value = 0.0202020202020202
$ws.Range("E7").Value = 0.005050505050505051
$ws.Range("F7").Value = 0.101010101010101
$ws.Range("J7").Value = 0.06565656565656566
$ws.Range("O7").Value = 0.0202020202020202
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.0707070707070707
$ws.Range("S7").Value = 0.4292929292929293
$ws.Range("B8").Value = 0.1029082774049217
$ws.Range("D8").Value = 0.01565995525727069
$ws.Range("F8").Value = 0.06263982102908278
$ws.Range("J8").Value = 0.1096196868008949
$ws.Range("O8").Value = 0.01565995525727069
$ws.Range("Q8").Value = 0.1655480984340045
$ws.Range("R8").Value = 0.116331096196868
$ws.Range("S8").Value = 0.4116331096196868
$ws.Range("B9").Value = 0.055
$ws.Range("D9").Value = 0.02
$ws.Range("F9").Value = 0.06
$ws.Range("J9").Value = 0.145
$ws.Range("O9").Value = 0.03
$ws.Range("Q9").Value = 0.17
$ws.Range("R9").Value = 0.075
$ws.Range("S9").Value = 0.445
$ws.Range("B10").Value = 0.09060955518945635
$ws.Range("D10").Value = 0.02059308072487644
$ws.Range("F10").Value = 0.06919275123558484
$ws.Range("J10").Value = 0.1350906095551895
$ws.Range("O10").Value = 0.0156507413509061
$ws.Range("Q10").Value = 0.186161449752883
$ws.Range("R10").Value = 0.1029654036243822
$ws.Range("S10").Value = 0.3797364085667216
$ws.Range("G11").Value = 0.1216216216216216
$ws.Range("J11").Value = 0.07094594594594594
$ws.Range("K11").Value = 0.1689189189189189
$ws.Range("L11").Value = 0.6148648648648649
$ws.Range("S11").Value = 0.02364864864864865
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1956521739130435
$ws.Range("K12").Value = 0.005434782608695652
$ws.Range("L12").Value = 0.02173913043478261
$ws.Range("S12").Value = 0.02717391304347826
$ws.Range("G13").Value = 0.7083333333333334
$ws.Range("S13").Value = 0.04166666666666666
$ws.Range("F15").Value = 0.01923076923076923
$ws.Range("H15").Value = 0.1201923076923077
$ws.Range("I15").Value = 0.0625
$ws.Range("J15").Value = 0.3701923076923077
$ws.Range("K15").Value = 0.07211538461538461
$ws.Range("M15").Value = 0.009615384615384616
$ws.Range("O15").Value = 0.1153846153846154
$ws.Range("S15").Value = 0.2307692307692308
$ws.Range("F16").Value = 0.01694915254237288
$ws.Range("H16").Value = 0.2259887005649718
$ws.Range("I16").Value = 0.0847457627118644
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.1129943502824859
$ws.Range("M16").Value = 0.01129943502824859
$ws.Range("O16").Value = 0.03389830508474576
$ws.Range("S16").Value = 0.1807909604519774
$ws.Range("F17").Value = 0.0202020202020202
$ws.Range("H17").Value = 0.1893939393939394
$ws.Range("I17").Value = 0.1262626262626263
$ws.Range("J17").Value = 0.3914141414141414
$ws.Range("K17").Value = 0.101010101010101
$ws.Range("M17").Value = 0.005050505050505051
$ws.Range("O17").Value = 0.06565656565656566
$ws.Range("S17").Value = 0.101010101010101
$ws.Range("F18").Value = 0.008928571428571428
$ws.Range("H18").Value = 0.1696428571428572
$ws.Range("I18").Value = 0.05357142857142857
$ws.Range("J18").Value = 0.4196428571428572
$ws.Range("K18").Value = 0.1071428571428571
$ws.Range("M18").Value = 0.03125
$ws.Range("N18").Value = 0.004464285714285714
$ws.Range("O18").Value = 0.04910714285714286
$ws.Range("S18").Value = 0.15625
$ws.Range("F19").Value = 0.01675552170601675
$ws.Range("H19").Value = 0.2086824067022087
$ws.Range("I19").Value = 0.08453922315308454
$ws.Range("J19").Value = 0.3465346534653465
$ws.Range("K19").Value = 0.1081492764661082
$ws.Range("M19").Value = 0.02817974105102818
$ws.Range("N19").Value = 0.001523229246001523
$ws.Range("O19").Value = 0.05864432597105865
$ws.Range("S19").Value = 0.146991622239147

Write-Host "Applied 105 cell updates"
